$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: G3 and H3 -> 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: D4 and E4 -> 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: D5 and E5 -> 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: D6 and E6 -> 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: H7 -> 1
$ws.Range("H7").Value = 1

# Row 8: H8 -> 1
$ws.Range("H8").Value = 1

# Row 9: D9 and E9 -> 1
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: D10 and E10 -> 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: D11 and E11 -> 1
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: D12 and E12 -> 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: H13 -> 1
$ws.Range("H13").Value = 1

# Row 14: D14 and E14 -> 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

# Row 15: D15 and E15 -> 1
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

# Row 16: H16 -> 1
$ws.Range("H16").Value = 1

# Row 17: D17 and E17 -> 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18: H18 -> 1
$ws.Range("H18").Value = 1
